$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.576.53"
$ws.Range("E2").Value = "  -0.27%  "
$ws.Range("D3").Value = "1.814.50"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +8.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.302"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0702"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0967"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.65%  "
$ws.Range("D12").Value = "2.074.23"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "1.823.35"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.657"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").Value = "34.539.00"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "0.0₃0801"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.13%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.123"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.38%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.24%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0533"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").Value = "1.403.97"
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -6.40%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.973"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.64%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("D48").Value = "1.975.15"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.36%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").Value = "0.0₆0128"
$ws.Range("E51").Value = "  -2.08%  "
